$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 46
$ws.Range("C37").Value = 9
$ws.Range("D37").Value = 16
$ws.Range("E37").Value = 20
$ws.Range("F37").Value = 71
$ws.Range("G37").Value = 91
